# Updates the cryptos list per the scraped data refresh (GitHub Actions run).
# For each changed cell we force a text NumberFormat before assigning the
# new value, then ClearFormats() to drop the temporary style again -- this
# keeps numeric-looking strings (e.g. "243.02", "0.07220") stored as text
# (matching the source data) instead of being auto-coerced into numbers by
# the normal Excel "smart" input parsing, while leaving cell styling as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "29.967.45"
Set-TextValue "E2" "  +0.33%  "
# Row 3
Set-TextValue "D3" "1.883.51"
Set-TextValue "E3" "  -0.23%  "
# Row 4
Set-TextValue "D4" "0.9995"
Set-TextValue "E4" "  -0.15%  "
# Row 5
Set-TextValue "D5" "0.7451"
Set-TextValue "E5" "  -2.96%  "
# Row 6
Set-TextValue "D6" "243.02"
Set-TextValue "E6" "  +0.19%  "
# Row 7
Set-TextValue "D7" "0.9996"
Set-TextValue "E7" "  -0.12%  "
# Row 8
Set-TextValue "D8" "0.3162"
Set-TextValue "E8" "  +0.97%  "
# Row 9
Set-TextValue "D9" "0.07220"
Set-TextValue "E9" "  +1.42%  "
# Row 10
Set-TextValue "E10" "  -2.74%  "
# Row 11
Set-TextValue "D11" "0.08343"
Set-TextValue "E11" "  -2.24%  "
# Row 12
Set-TextValue "B12" "WrappedEther"
Set-TextValue "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.948.69"
Set-TextValue "E12" "  +2.11%  "
# Row 13
Set-TextValue "B13" "Polygon"
Set-TextValue "C13" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D13" "0.7558"
Set-TextValue "E13" "  -0.92%  "
# Row 14
Set-TextValue "B14" "Polkadot"
Set-TextValue "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.421"
Set-TextValue "E14" "  +1.04%  "
# Row 15
Set-TextValue "D15" "92.54"
Set-TextValue "E15" "  -1.08%  "
# Row 16
Set-TextValue "D16" "6.151"
Set-TextValue "E16" "  +0.23%  "
# Row 17
Set-TextValue "D17" "30.043.22"
Set-TextValue "E17" "  +0.42%  "
# Row 18
Set-TextValue "D18" "249.82"
Set-TextValue "E18" "  +2.37%  "
# Row 19
Set-TextValue "D19" "13.61"
Set-TextValue "E19" "  -0.97%  "
# Row 20
Set-TextValue "D20" "0.000007867"
Set-TextValue "E20" "  +0.64%  "
# Row 21
Set-TextValue "D21" "2.204.03"
Set-TextValue "E21" "  +2.20%  "
# Row 22
Set-TextValue "E22" "  +0.11%  "
# Row 23
Set-TextValue "D23" "7.999"
Set-TextValue "E23" "  -0.11%  "
# Row 24
Set-TextValue "D24" "0.9990"
Set-TextValue "E24" "  -0.18%  "
# Row 25
Set-TextValue "D25" "0.1565"
Set-TextValue "E25" "  -3.77%  "
# Row 26
Set-TextValue "D26" "9.293"
Set-TextValue "E26" "  -1.02%  "
# Row 27
Set-TextValue "D27" "165.48"
Set-TextValue "E27" "  +1.40%  "
# Row 28
Set-TextValue "E28" "  -0.41%  "
# Row 29
Set-TextValue "D29" "2.040"
Set-TextValue "E29" "  +0.22%  "
# Row 30
Set-TextValue "D30" "1.484"
Set-TextValue "E30" "  -2.24%  "
# Row 31
Set-TextValue "E31" "  +2.26%  "
# Row 32
Set-TextValue "E32" "  -0.04%  "
# Row 33
Set-TextValue "D33" "4.231"
Set-TextValue "E33" "  +2.54%  "
# Row 34
Set-TextValue "D34" "0.05371"
Set-TextValue "E34" "  -1.30%  "
# Row 35
Set-TextValue "E35" "  +0.98%  "
# Row 36
Set-TextValue "D36" "0.7559"
Set-TextValue "E36" "  +1.42%  "
# Row 37
Set-TextValue "D37" "0.9936"
Set-TextValue "E37" "  -0.70%  "
# Row 38
Set-TextValue "D38" "2.706"
Set-TextValue "E38" "  -0.02%  "
# Row 39
Set-TextValue "D39" "0.01966"
Set-TextValue "E39" "  +0.96%  "
# Row 40
Set-TextValue "E40" "  -0.63%  "
# Row 41
Set-TextValue "D41" "0.4562"
Set-TextValue "E41" "  +2.07%  "
# Row 42
Set-TextValue "D42" "1.106.85"
Set-TextValue "E42" "  +0.54%  "
# Row 43
Set-TextValue "B43" "Aave"
Set-TextValue "C43" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "72.76"
# Row 44
Set-TextValue "B44" "FraxShare"
Set-TextValue "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "6.051"
Set-TextValue "E44" "  -0.47%  "
# Row 45
Set-TextValue "D45" "0.8693"
Set-TextValue "E45" "  +1.46%  "
# Row 46
Set-TextValue "D46" "104.44"
Set-TextValue "E46" "  +1.41%  "
# Row 47
Set-TextValue "D47" "1.001"
Set-TextValue "E47" "  +0.01%  "
# Row 48
Set-TextValue "D48" "1.866"
Set-TextValue "E48" "  -0.21%  "
# Row 49
Set-TextValue "D49" "7.614"
Set-TextValue "E49" "  -0.67%  "
# Row 50
Set-TextValue "D50" "2.093.42"
Set-TextValue "E50" "  +2.67%  "
# Row 51
Set-TextValue "D51" "9.546"
Set-TextValue "E51" "  -1.87%  "
